# Update build version timestamp strings across the workbook.
$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $val = $cell.Value2
        if ($null -ne $val -and $val -is [string] -and $val.Contains($oldStamp)) {
            $cell.Value2 = $val.Replace($oldStamp, $newStamp)
        }
    }
}
